$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for price cells whose new value looks like a number,
# so Excel doesn't silently convert the literal digit-string into a Double.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = '22.412.14'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '1.570.86'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '289.96'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.3745'
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("D8").Value = '49.38'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '0.3363'
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").Value = '0.07457'
$ws.Range("E10").Value = '  -2.88%  '
$ws.Range("D11").Value = '1.127'
$ws.Range("E11").Value = '  -3.40%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '20.94'
$ws.Range("E13").Value = '  -2.12%  '
$ws.Range("D14").Value = '5.912'
$ws.Range("E14").Value = '  -1.92%  '
$ws.Range("D15").Value = '6.860'
$ws.Range("E15").Value = '  -1.23%  '
$ws.Range("D16").Value = '1.568.88'
$ws.Range("E16").Value = '  -0.93%  '
$ws.Range("D17").Value = '0.00001115'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '89.12'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").Value = '0.06695'
$ws.Range("E19").Value = '  -0.39%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").Value = '6.160'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").Value = '16.16'
$ws.Range("E22").Value = '  -2.97%  '
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("D24").Value = '22.402.08'
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '2.366'
$ws.Range("E25").Value = '  -1.03%  '
$ws.Range("D26").Value = '2.540'
$ws.Range("E26").Value = '  -9.04%  '
$ws.Range("D27").Value = '19.99'
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("D28").Value = '146.95'
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").Value = '4.997'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").Value = '124.59'
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("D31").Value = '1.743.35'
$ws.Range("E31").Value = '  -0.57%  '
$ws.Range("D32").Value = '0.9934'
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("D33").Value = '1.957'
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("D34").Value = '5.898'
$ws.Range("E34").Value = '  -5.34%  '
$ws.Range("D35").Value = '9.711'
$ws.Range("E35").Value = '  -4.36%  '
$ws.Range("D36").Value = '0.08384'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '1.372'
$ws.Range("E37").Value = '  +3.29%  '
$ws.Range("D38").Value = '0.02448'
$ws.Range("E38").Value = '  -3.80%  '
$ws.Range("D39").Value = '0.06462'
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("D40").Value = '0.2244'
$ws.Range("E40").Value = '  -3.50%  '
$ws.Range("D41").Value = '5.367'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("D42").Value = '11.31'
$ws.Range("E42").Value = '  -3.77%  '
$ws.Range("D43").Value = '0.6202'
$ws.Range("E43").Value = '  -3.30%  '
$ws.Range("D44").Value = '1.003'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").Value = '14.04'
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").Value = '3.809'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("D47").Value = '0.5786'
$ws.Range("E47").Value = '  -3.55%  '
$ws.Range("D48").Value = '2.055'
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").Value = '125.38'
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("E50").Value = '  -4.16%  '
$ws.Range("E51").Value = '  +0.05%  '
